$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.535459160804749
$ws.Range("B1").Value = 1.686869978904724
$ws.Range("C1").Value = 2.003501176834106
$ws.Range("D1").Value = 3.066848278045654
$ws.Range("E1").Value = 1.365401148796082
